# Updated cryptos list (price + 1h volume refresh), matching the
# GitHub Actions scheduled scrape commit.
#
# Note: several "Price" cells hold values that look numeric (e.g. "572.50",
# "1.00", "0.0000117") but must stay plain text, matching the source sheet's
# inline-string cells. Setting NumberFormat to "@" (Text) before assigning
# the value keeps the COM layer from auto-coercing them to numbers/losing
# trailing zeros, and resetting Style back to "Normal" afterwards avoids
# leaving a stray text-format style on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.527.37"
$ws.Range("D3").Value = "3.331.93"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.617"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.70%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -3.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.401"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.20%  "
$ws.Range("D12").Value = "3.911.75"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.13%  "
$ws.Range("D15").Value = "66.623.13"
$ws.Range("E15").Value = "  -4.09%  "
$ws.Range("E16").Value = "  -2.50%  "
$ws.Range("D17").Value = "3.317.71"
$ws.Range("E17").Value = "  -2.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "437.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.517"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000117"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.190"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("E31").Value = "  -6.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "162.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "27.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.05%  "
$ws.Range("D39").Value = "2.809.01"
$ws.Range("E39").Value = "  +2.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.790"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.23%  "
$ws.Range("E41").Value = "  -4.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0665"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "320.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.03%  "
$ws.Range("E48").Value = "  -3.97%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.102"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.975"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.14%  "
